# erillishaku_oidilla: support excel import with a person-oid, and allow
# several applicants to share the same birth date.
#
# Row 2 used to carry a Finnish "henkilötunnus" (national id, C2) together
# with a birth date (D2). We now clear the national id, keep the birth date
# (re-typed as plain text so Excel/readers don't reinterpret it as a date
# serial), and introduce a new "Hakija-oid" (applicant oid) value in E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Syntymäaika (D2) must stay textual, so force a text number format before
# writing the value - otherwise "1.1.1980" gets auto-parsed into a date
# serial number.
$dateCell = $ws.Range("D2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1.1.1980"

# Henkilötunnus (C2) is no longer collected - clear it out.
$ws.Range("C2").Value = ""

# Hakija-oid (E2) is new.
$ws.Range("E2").Value = "hakija1"

# Reflect where the user last clicked while editing.
$ws.Range("D3").Select() | Out-Null
